$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.61 = 18179.77 pesos`n✅ 18179.77 pesos = 4.58 = 957.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 217
$wsTasas.Range("O10").Value = 3945.01
$wsTasas.Range("N12").Value = 3970
$wsTasas.Range("O12").Value = 209
